# "Moved new docs to gh-pages branch" — append the two meeting-log entries that
# were logged after the prior save: fill in the (already-present, blank) table
# row 25 and append a brand-new row 26, then keep the table / selection in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meeting Logs")
$lo = $ws.ListObjects.Item("Table2")

# Grow the table by one row (row 25 already exists as a blank table row, so
# this produces the new row 26 and extends Table2 to A1:D26).
$lo.ListRows.Add() | Out-Null

# --- Row 25 (Fri 11/30/18 @ Leep2) ---------------------------------------
$ws.Range("A25").Value = 43434
$ws.Range("B25").Value = "Leep2"
$ws.Range("C25").Value = "Gage, Jacob"
$ws.Range("D25").Value = "Worked out artifacts and documentation that still needed to be done, set up another meeting to work on these"

# --- Row 26 (Sun 12/2/18 @ Spahr Library) --------------------------------
# Entered in the same order the workbook's shared-string table shows it was
# originally typed: date, then Members, then Things Discussed, then Location.
$ws.Range("A26").Value = 43436
$ws.Range("A26").NumberFormat = "m/d/yy"
$ws.Range("C26").Value = "Gage, Jacob, Zach"
$ws.Range("D26").Value = "Working on remaining artifacts and documentation"
$ws.Range("B26").Value = "Spahr Library"

# The "Things Discussed" column wraps text for every data row.
$ws.Range("D25:D26").WrapText = $true

# One existing row's Members cell also picked up a wrap-text format.
$ws.Range("C13").WrapText = $true

# Leave the cursor where the author left it.
$ws.Range("E2").Select() | Out-Null
